$p = $ppt.ActivePresentation
$d1 = $p.Designs.Add()
$d2 = $p.Designs.Add()
